$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: quantity 98 -> 95, title "Backpa" -> "Backpack"
$ws.Range("C2").Value = 95
$ws.Range("D2").Value = "Backpack"

# Row 3: quantity 8 -> 6
$ws.Range("C3").Value = 6

# Row 4: quantity 9 -> 7
$ws.Range("C4").Value = 7
